$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6: A6 = "Timmy Hellfire"
$ws.Range("A6").Value = "Timmy Hellfire"

# New row 8: C8 = "purple"
$ws.Range("C8").Value = "purple"

# A4 changes from "Prat-a" to "Viveeon Eastwood"
$ws.Range("A4").Value = "Viveeon Eastwood"

# Update selection to A5
$ws.Range("A5").Select()
